# Updated Chapter 2 and 3 to PeptideShaker v0.23.0 and SearchGUI v1.15.0.
#
# 1) Expand "... folder. You should see the following:" into
#    "... folder. When done loading you should see the following:"
# 2) The extra sentence pushes the following image onto the next page,
#    so the footer's cached PAGE field result flips from 2 -> 3.

$d = $word.ActiveDocument

# --- 1. Body text: insert "When done loading " before "you should see the following:" ---
$body = $d.Content
$found = $body.Find.Execute(
    " folder. You should see the following:",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    " folder. When done loading you should see the following:",
    2)

if (-not $found) {
    throw "Could not find the target sentence to update."
}

# --- 2. Footer: cached PAGE field result 2 -> 3 (re-pagination side effect) ---
$sec = $d.Sections(1)
$ftr = $sec.Footers(1)
$frange = $ftr.Range
$frange.Find.Execute(
    "2",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "3",
    2)
